# "remove column from alcohol data"
#
# The data sheet (Sheet1) had a duplicate/trailing measurement column (N)
# that needs to go away. Deleting column M shifts the old column N left
# into the M position (and everything else stays put), which is exactly
# the A1:N119 -> A1:M119 shrink seen in the target workbook.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# Remove the column; remaining columns to the right (just N) shift left.
$ws1.Columns("M:M").Delete()

# Reflect the author's resulting view state: selection parked on the new
# last column, zoomed in from the very-zoomed-out 55% to 95%.
$ws1.Range("M1").Select() | Out-Null
$excel.ActiveWindow.Zoom = 95

# The other (empty) sheets in the workbook end up re-zoomed to 95% too.
$ws2.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 95

$ws3.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 95

# Leave the original sheet active/selected, matching tabSelected="true".
$ws1.Activate() | Out-Null
